$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2986.1738
$ws.Range("I64").Value = 2945.111
$ws.Range("J64").Value = 3134
$ws.Range("K64").Value = 2945.111
$ws.Range("L64").Value = 3134
$ws.Range("M64").Value = -2697.111
$ws.Range("N64").Value = -3630
$ws.Range("H67").Value = 2986.1738
$ws.Range("I67").Value = 2945.111
$ws.Range("J67").Value = 3134
$ws.Range("K67").Value = 2945.111
$ws.Range("L67").Value = 3134
$ws.Range("M67").Value = -2087.111
$ws.Range("N67").Value = -4850
$ws.Range("H97").Value = 1800
$ws.Range("J97").Value = 2000
$ws.Range("L97").Value = 6000
$ws.Range("N97").Value = -6992
$ws.Range("H98").Value = 389537.94
$ws.Range("I98").Value = 6486.8667
$ws.Range("J98").Value = 911880.25
$ws.Range("K98").Value = 6486.8667
$ws.Range("L98").Value = 911880.25
$ws.Range("M98").Value = -4988.8667
$ws.Range("N98").Value = -914876.25
$ws.Range("H100").Value = 1974.3478
$ws.Range("I100").Value = 1134.6666
$ws.Range("J100").Value = 4997.2
$ws.Range("K100").Value = 1134.6666
$ws.Range("L100").Value = 4997.2
$ws.Range("M100").Value = -593.6666
$ws.Range("N100").Value = -6079.2
$ws.Range("H103").Value = 62500720
$ws.Range("I103").Value = 893.5
$ws.Range("J103").Value = 125000550
$ws.Range("K103").Value = 2680.5
$ws.Range("L103").Value = 375001650
$ws.Range("M103").Value = -2094.5
$ws.Range("N103").Value = -375002822
$ws.Range("H106").Value = 3180
$ws.Range("I106").Value = 3007.1428
$ws.Range("J106").Value = 5600
$ws.Range("K106").Value = 3007.1428
$ws.Range("L106").Value = 5600
$ws.Range("M106").Value = -2376.1428
$ws.Range("N106").Value = -6862
$ws.Range("H107").Value = 967.05884
$ws.Range("I107").Value = 1296
$ws.Range("J107").Value = 497.14285
$ws.Range("K107").Value = 1296
$ws.Range("L107").Value = 497.14285
$ws.Range("M107").Value = 624
$ws.Range("N107").Value = -4337.14285
$ws.Range("H109").Value = 18460.5
$ws.Range("J109").Value = 18460.5
$ws.Range("L109").Value = 18460.5
$ws.Range("N109").Value = -21234.5
$ws.Range("H111").Value = 1206.9565
$ws.Range("I111").Value = 993.06665
$ws.Range("J111").Value = 1608
$ws.Range("K111").Value = 2979.19995
$ws.Range("L111").Value = 4824
$ws.Range("M111").Value = 87.80004999999983
$ws.Range("N111").Value = -10958
$ws.Range("H112").Value = 5740.543
$ws.Range("J112").Value = 6560.6333
$ws.Range("L112").Value = 19681.8999
$ws.Range("N112").Value = -21897.8999
$ws.Range("H115").Value = 633.5
$ws.Range("I115").Value = 412.14285
$ws.Range("J115").Value = 1150
$ws.Range("K115").Value = 1236.42855
$ws.Range("L115").Value = 3450
$ws.Range("M115").Value = 330.5714499999999
$ws.Range("N115").Value = -6584
$ws.Range("H118").Value = 1246
$ws.Range("I118").Value = 232
$ws.Range("J118").Value = 2260
$ws.Range("K118").Value = 696
$ws.Range("L118").Value = 6780
$ws.Range("M118").Value = 961
$ws.Range("N118").Value = -10094
$ws.Range("H122").Value = 389537.94
$ws.Range("I122").Value = 6486.8667
$ws.Range("J122").Value = 911880.25
$ws.Range("K122").Value = 19460.6001
$ws.Range("L122").Value = 2735640.75
$ws.Range("M122").Value = -17010.6001
$ws.Range("N122").Value = -2740540.75
$ws.Range("H124").Value = 45000
$ws.Range("J124").Value = 45000
$ws.Range("L124").Value = 45000
$ws.Range("N124").Value = -54820
$ws.Range("H126").Value = 49980
$ws.Range("J126").Value = 49980
$ws.Range("L126").Value = 49980
$ws.Range("N126").Value = -59860
$ws.Range("H141").Value = 1494.25
$ws.Range("I141").Value = 1522
$ws.Range("J141").Value = 1355.5
$ws.Range("K141").Value = 4566
$ws.Range("L141").Value = 4066.5
$ws.Range("M141").Value = 614
$ws.Range("N141").Value = -14426.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1879.2826
$ws.Range("I132").Value = 1096.9656
$ws.Range("K132").Value = 3290.8968
$ws.Range("M132").Value = -760.8968

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1864.5883
$ws.Range("I58").Value = 1899.2
$ws.Range("J58").Value = 1815.1428
$ws.Range("K58").Value = 1899.2
$ws.Range("L58").Value = 1815.1428
$ws.Range("M58").Value = -1696.2
$ws.Range("N58").Value = -2221.1428
$ws.Range("H136").Value = 1864.5883
$ws.Range("I136").Value = 1899.2
$ws.Range("J136").Value = 1815.1428
$ws.Range("K136").Value = 5697.6
$ws.Range("L136").Value = 5445.428400000001
$ws.Range("M136").Value = -3147.6
$ws.Range("N136").Value = -10545.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 962.62964
$ws.Range("I68").Value = 950
$ws.Range("J68").Value = 970.05884
$ws.Range("K68").Value = 2850
$ws.Range("L68").Value = 2910.17652
$ws.Range("M68").Value = -2039
$ws.Range("N68").Value = -4532.17652
$ws.Range("H71").Value = 962.62964
$ws.Range("I71").Value = 950
$ws.Range("J71").Value = 970.05884
$ws.Range("K71").Value = 8550
$ws.Range("L71").Value = 8730.529560000001
$ws.Range("M71").Value = -4494
$ws.Range("N71").Value = -16842.52956
$ws.Range("H131").Value = 9260141
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 9260141
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 27780423
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -27790503
$ws.Range("H132").Value = 1730
$ws.Range("I132").Value = 1850
$ws.Range("J132").Value = 1592.8572
$ws.Range("K132").Value = 16650
$ws.Range("L132").Value = 14335.7148
$ws.Range("M132").Value = -14120
$ws.Range("N132").Value = -19395.7148

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2566.5715
$ws.Range("I7").Value = 2329.3333
$ws.Range("J7").Value = 3990
$ws.Range("K7").Value = 2329.3333
$ws.Range("L7").Value = 3990
$ws.Range("M7").Value = -2217.3333
$ws.Range("N7").Value = -4214
$ws.Range("H40").Value = 1718.8
$ws.Range("I40").Value = 1522.2142
$ws.Range("J40").Value = 2177.5
$ws.Range("K40").Value = 1522.2142
$ws.Range("L40").Value = 2177.5
$ws.Range("M40").Value = -1386.2142
$ws.Range("N40").Value = -2449.5
$ws.Range("H122").Value = 2642.8572
$ws.Range("I122").Value = 2670
$ws.Range("J122").Value = 2575
$ws.Range("K122").Value = 8010
$ws.Range("L122").Value = 7725
$ws.Range("M122").Value = -5560
$ws.Range("N122").Value = -12625
$ws.Range("H126").Value = 2566.5715
$ws.Range("I126").Value = 2329.3333
$ws.Range("J126").Value = 3990
$ws.Range("K126").Value = 6987.999899999999
$ws.Range("L126").Value = 11970
$ws.Range("M126").Value = -4517.999899999999
$ws.Range("N126").Value = -16910
$ws.Range("H136").Value = 4285.1304
$ws.Range("I136").Value = 3237.0908
$ws.Range("J136").Value = 5245.8335
$ws.Range("K136").Value = 9711.2724
$ws.Range("L136").Value = 15737.5005
$ws.Range("M136").Value = -7161.2724
$ws.Range("N136").Value = -20837.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1336.069
$ws.Range("I122").Value = 1179.3636
$ws.Range("J122").Value = 1828.5714
$ws.Range("K122").Value = 3538.0908
$ws.Range("L122").Value = 5485.7142
$ws.Range("M122").Value = -1088.0908
$ws.Range("N122").Value = -10385.7142
$ws.Range("H126").Value = 1543.7097
$ws.Range("I126").Value = 1016.4091
$ws.Range("J126").Value = 2832.6667
$ws.Range("K126").Value = 3049.2273
$ws.Range("L126").Value = 8498.000100000001
$ws.Range("M126").Value = -579.2273
$ws.Range("N126").Value = -13438.0001
